# Generate Report for Handoff
# This script swaps the "a7042400..." and "6cc6d73e..." rows (the report was
# regenerated and the two source files traded places in the sort order), and
# updates the 6cc6d73e row's status from "Handed back: in sync with en-US" to
# "Ready for handoff" (including the new handoff timestamps and the stale
# handback warning that now appears for it).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

# Row 2 becomes the a7042400 file (was 6cc6d73e)
$ov.Range("A2").Value = "a7042400-8349-4ef4-9d00-803d03e8f618.md"
$ov.Range("B2").Value = "e2e\a7042400-8349-4ef4-9d00-803d03e8f618.md"

# Row 3 becomes the 6cc6d73e file (was a7042400) and is now "Ready for handoff"
$ov.Range("A3").Value = "6cc6d73e-a009-46b1-b8fd-922c4b59b5c0.md"
$ov.Range("B3").Value = "e2e\6cc6d73e-a009-46b1-b8fd-922c4b59b5c0.md"
$ov.Range("E3").Value = "Ready for handoff"
$ov.Range("F3").Value = "Ready for handoff"
$ov.Range("G3").Value = "2016-08-15 20:45:05"

# The B2/B3 hyperlinks keep pointing at the same target URLs (rId2 -> 6cc6d73e
# URL, rId3 -> a7042400 URL) but their displayed text now matches the new
# cell text above, so rebuild them against the same addresses with swapped
# display text.
$ov.Hyperlinks.Delete()
$ov.Hyperlinks.Add($ov.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2b04e7a825218e3d74ff95af77cea65a52aa8612/e2e/6cc6d73e-a009-46b1-b8fd-922c4b59b5c0.md", "", "", "e2e\a7042400-8349-4ef4-9d00-803d03e8f618.md")
$ov.Hyperlinks.Add($ov.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2b04e7a825218e3d74ff95af77cea65a52aa8612/e2e/a7042400-8349-4ef4-9d00-803d03e8f618.md", "", "", "e2e\6cc6d73e-a009-46b1-b8fd-922c4b59b5c0.md")

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

# Row 2 becomes the a7042400 file (was 6cc6d73e) - still "Handed back" status
$zh.Range("A2").Value = "a7042400-8349-4ef4-9d00-803d03e8f618.md"
$zh.Range("G2").Value = "a7042400-8349-4ef4-9d00-803d03e8f618.ac66859d14d4501f31d5e816897883c12bcfcd2e.zh-cn.xlf"
$zh.Range("I2").Value = "a7042400-8349-4ef4-9d00-803d03e8f618.md"
$zh.Range("J2").Value = "a7042400-8349-4ef4-9d00-803d03e8f618.ac66859d14d4501f31d5e816897883c12bcfcd2e.zh-cn.xlf"

# Row 3 becomes the 6cc6d73e file (was a7042400) - now "Ready for handoff"
# with a fresh handoff file/time and a stale-handback error message.
$zh.Range("A3").Value = "6cc6d73e-a009-46b1-b8fd-922c4b59b5c0.md"
$zh.Range("C3").Value = "Ready for handoff"
$zh.Range("G3").Value = "6cc6d73e-a009-46b1-b8fd-922c4b59b5c0.684d8844e0884ae608929bad0eabacf861d159b2.zh-cn.xlf"
$zh.Range("H3").Value = "2016-08-15 20:44:56"
$zh.Range("I3").Value = "6cc6d73e-a009-46b1-b8fd-922c4b59b5c0.md"
$zh.Range("J3").Value = "6cc6d73e-a009-46b1-b8fd-922c4b59b5c0.684d8844e0884ae608929bad0eabacf861d159b2.zh-cn.xlf"
$zh.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2b04e7a825218e3d74ff95af77cea65a52aa8612/e2e/6cc6d73e-a009-46b1-b8fd-922c4b59b5c0.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/aa95900b6240febcd94ab35a36b897ba01241851/e2e/6cc6d73e-a009-46b1-b8fd-922c4b59b5c0.md."

# The "Error Detail" column is now wide enough to show the message above.
$zh.Columns.Item(16).ColumnWidth = 39.14

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

# Row 2 becomes the a7042400 file (was 6cc6d73e) - still "Handed back" status
$de.Range("A2").Value = "a7042400-8349-4ef4-9d00-803d03e8f618.md"
$de.Range("G2").Value = "a7042400-8349-4ef4-9d00-803d03e8f618.ac66859d14d4501f31d5e816897883c12bcfcd2e.de-de.xlf"
$de.Range("I2").Value = "a7042400-8349-4ef4-9d00-803d03e8f618.md"
$de.Range("J2").Value = "a7042400-8349-4ef4-9d00-803d03e8f618.ac66859d14d4501f31d5e816897883c12bcfcd2e.de-de.xlf"

# Row 3 becomes the 6cc6d73e file (was a7042400) - now "Ready for handoff"
# with a fresh handoff file/time and a stale-handback error message.
$de.Range("A3").Value = "6cc6d73e-a009-46b1-b8fd-922c4b59b5c0.md"
$de.Range("C3").Value = "Ready for handoff"
$de.Range("G3").Value = "6cc6d73e-a009-46b1-b8fd-922c4b59b5c0.684d8844e0884ae608929bad0eabacf861d159b2.de-de.xlf"
$de.Range("H3").Value = "2016-08-15 20:45:05"
$de.Range("I3").Value = "6cc6d73e-a009-46b1-b8fd-922c4b59b5c0.md"
$de.Range("J3").Value = "6cc6d73e-a009-46b1-b8fd-922c4b59b5c0.684d8844e0884ae608929bad0eabacf861d159b2.de-de.xlf"
$de.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2b04e7a825218e3d74ff95af77cea65a52aa8612/e2e/6cc6d73e-a009-46b1-b8fd-922c4b59b5c0.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/aa95900b6240febcd94ab35a36b897ba01241851/e2e/6cc6d73e-a009-46b1-b8fd-922c4b59b5c0.md."

# The "Error Detail" column is now wide enough to show the message above.
$de.Columns.Item(16).ColumnWidth = 39.14
